# Pins.xlsx update: rebuild the pin table with new columns (In/Out, Analog/Digital,
# Motor Driver) and updated/added pin rows (encoders, corrected IR/Teensy numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale cells/values don't linger.
$ws.Cells.Clear()

# ---- Header row (bold) ----
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "Teensy"
$ws.Range("C1").Value = "Analog Name"
$ws.Range("D1").Value = "In/Out"
$ws.Range("E1").Value = "Analog/Digital"
$ws.Range("F1").Value = "Motor Driver"
$ws.Range("A1:F1").Font.Bold = $true

# ---- Motor driver pins ----
$ws.Range("A2").Value = "L Motor Speed"
$ws.Range("B2").Value = 7
$ws.Range("D2").Value = "Out"
$ws.Range("E2").Value = "Analog"
$ws.Range("F2").Value = "PWMA"

$ws.Range("A3").Value = "L Motor Forward"
$ws.Range("B3").Value = 8
$ws.Range("D3").Value = "Out"
$ws.Range("E3").Value = "Digital"
$ws.Range("F3").Value = "AIN2"

$ws.Range("A4").Value = "L Motor Reverse"
$ws.Range("B4").Value = 9
$ws.Range("D4").Value = "Out"
$ws.Range("E4").Value = "Digital"
$ws.Range("F4").Value = "AIN1"

$ws.Range("A5").Value = "R Motor Speed"
$ws.Range("B5").Value = 10
$ws.Range("D5").Value = "Out"
$ws.Range("E5").Value = "Analog"
$ws.Range("F5").Value = "PWMB"

$ws.Range("A6").Value = "R Motor Forward"
$ws.Range("B6").Value = 11
$ws.Range("D6").Value = "Out"
$ws.Range("E6").Value = "Digital"
$ws.Range("F6").Value = "BIN2"

$ws.Range("A7").Value = "R Motor Reverse"
$ws.Range("B7").Value = 12
$ws.Range("D7").Value = "Out"
$ws.Range("E7").Value = "Digital"
$ws.Range("F7").Value = "BIN1"

$ws.Range("A8").Value = "STBY"
$ws.Range("B8").Value = 13
$ws.Range("D8").Value = "Out"
$ws.Range("E8").Value = "Digital"
$ws.Range("F8").Value = "STBY"

# ---- IR sensor pins ----
$ws.Range("A9").Value = "IR Right"
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = "A13"
$ws.Range("D9").Value = "In"
$ws.Range("E9").Value = "Analog"

$ws.Range("A10").Value = "IR Middle"
$ws.Range("B10").Value = 33
$ws.Range("C10").Value = "A15"
$ws.Range("D10").Value = "In"
$ws.Range("E10").Value = "Analog"

$ws.Range("A11").Value = "IR Left"
$ws.Range("B11").Value = 34
$ws.Range("C11").Value = "A14"
$ws.Range("D11").Value = "In"
$ws.Range("E11").Value = "Analog"

# ---- Encoder pins (new) ----
$ws.Range("A12").Value = "L Encoder A"
$ws.Range("B12").Value = 38
$ws.Range("D12").Value = "In"
$ws.Range("E12").Value = "Digital"

$ws.Range("A13").Value = "L Encoder B"
$ws.Range("B13").Value = 37
$ws.Range("D13").Value = "In"
$ws.Range("E13").Value = "Digital"

$ws.Range("A14").Value = "R Encoder A"
$ws.Range("B14").Value = 27
$ws.Range("D14").Value = "In"
$ws.Range("E14").Value = "Digital"

$ws.Range("A15").Value = "R Encoder B"
$ws.Range("B15").Value = 28
$ws.Range("D15").Value = "In"
$ws.Range("E15").Value = "Digital"

# ---- Column widths (best-fit to content, like the original column A) ----
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# ---- Page setup / selection to mirror the saved view state ----
$ws.PageSetup.Orientation = 1
$ws.Range("F12").Select() | Out-Null
